$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.554.71'
$ws.Range("E2").Value = '  -2.18%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.881.97'
$ws.Range("E3").Value = '  -2.27%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.34'
$ws.Range("E5").Value = '  -4.34%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.06'
$ws.Range("E6").Value = '  -3.46%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.500'
$ws.Range("E8").Value = '  -0.98%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.880.59'
$ws.Range("E9").Value = '  -2.26%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.84'
$ws.Range("E10").Value = '  -2.45%  '

# Row 11
$ws.Range("E11").Value = '  -2.28%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.428'
$ws.Range("E12").Value = '  -1.99%  '

# Row 13
$ws.Range("E13").Value = '  -1.70%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.54'
$ws.Range("E14").Value = '  -2.73%  '

# Row 15
$ws.Range("E15").Value = '  -0.22%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.360.74'

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.527.32'
$ws.Range("E17").Value = '  -2.19%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.894.20'
$ws.Range("E18").Value = '  -2.00%  '

# Row 19
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.47'
$ws.Range("E19").Value = '  -2.92%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '428.56'
$ws.Range("E20").Value = '  -2.04%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.99'
$ws.Range("E21").Value = '  -2.92%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.649'
$ws.Range("E22").Value = '  -2.18%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.78'
$ws.Range("E23").Value = '  -3.27%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.74'
$ws.Range("E24").Value = '  -2.33%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.86'
$ws.Range("E25").Value = '  +0.75%  '

# Row 26
$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.02'
$ws.Range("E26").Value = '  -10.87%  '

# Row 27
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.00%  '

# Row 28
$ws.Range("E28").Value = '  -5.65%  '

# Row 29
$ws.Range("E29").Value = '  +6.57%  '

# Row 30
$ws.Range("E30").Value = '  -3.79%  '

# Row 31
$ws.Range("E31").Value = '  -4.61%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.02'
$ws.Range("E32").Value = '  -9.34%  '

# Row 33
$ws.Range("E33").Value = '  +0.01%  '

# Row 34
$ws.Range("E34").Value = '  -2.13%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.32'
$ws.Range("E35").Value = '  -3.82%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.953'
$ws.Range("E36").Value = '  -3.89%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.34'
$ws.Range("E37").Value = '  -4.60%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.75'
$ws.Range("E38").Value = '  -1.72%  '

# Row 39
$ws.Range("E39").Value = '  -7.02%  '

# Row 40
$ws.Range("E40").Value = '  -6.02%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.13'
$ws.Range("E41").Value = '  -3.39%  '

# Row 42
$ws.Range("E42").Value = '  -3.75%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.09'
$ws.Range("E43").Value = '  -0.89%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.265'
$ws.Range("E44").Value = '  -5.02%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.676.53'
$ws.Range("E45").Value = '  -0.30%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '132.07'
$ws.Range("E46").Value = '  -2.64%  '

# Row 47
$ws.Range("E47").Value = '  -1.16%  '

# Row 48
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '342.38'
$ws.Range("E48").Value = '  -4.30%  '

# Row 49
$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  -0.01%  '

# Row 50
$ws.Range("E50").Value = '  -2.01%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.36'
$ws.Range("E51").Value = '  -5.53%  '
